# 3. ödev not güncellemesi
# Melike Yoğurtçu: 0 -> 60  (Ödev 3 / column E, row 14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the grade for Melike Yoğurtçu's 3rd homework (Ödev 3)
$ws.Range("E14").Value = 60

# Reflect the new active cell/selection recorded in the saved file
$ws.Range("E17").Select()
